$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 20 values (2025Q2)
$ws.Range("C20").Value = 348
$ws.Range("D20").Value = 270
$ws.Range("F20").Value = 83.59133126934985

# Add new row 21 (2025Q3)
$ws.Range("A21").Value = "2025Q3"
$ws.Range("B21").Value = "2025Q3"
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0.2873563218390804
